$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 48/49 swap: BOLO <-> CoinbaseStockToken order swapped, with updated price/volume ---
$ws.Range("B48").Value = "CoinbaseStockToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D48").Value = "'0.002280"
$ws.Range("E48").Value = "'-0.55%"

$ws.Range("B49").Value = "BOLO"
$ws.Range("C49").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D49").Value = "'0.003207"
$ws.Range("E49").Value = "'6.74%"

# --- Price / Volume(1h) updates for all other changed rows ---
$ws.Range("D2").Value = "'328.68"
$ws.Range("E2").Value = "'1.70%"
$ws.Range("D3").Value = "'41.44"
$ws.Range("E3").Value = "'4.81%"
$ws.Range("D4").Value = "'5.620"
$ws.Range("E4").Value = "'-2.36%"
$ws.Range("D5").Value = "'0.08174"
$ws.Range("E5").Value = "'2.05%"
$ws.Range("D6").Value = "'2.024"
$ws.Range("E6").Value = "'1.92%"
$ws.Range("D7").Value = "'8.738"
$ws.Range("E7").Value = "'1.29%"
$ws.Range("D8").Value = "'4.525"
$ws.Range("E8").Value = "'-0.54%"
$ws.Range("D9").Value = "'2.963"
$ws.Range("E9").Value = "'0.34%"
$ws.Range("D10").Value = "'0.9213"
$ws.Range("E10").Value = "'-0.87%"
$ws.Range("D11").Value = "'0.1278"
$ws.Range("E11").Value = "'0.65%"
$ws.Range("E12").Value = "'-0.08%"
$ws.Range("D13").Value = "'0.09367"
$ws.Range("E13").Value = "'2.73%"
$ws.Range("D14").Value = "'0.03811"
$ws.Range("E14").Value = "'3.20%"
$ws.Range("D15").Value = "'0.1060"
$ws.Range("E15").Value = "'1.13%"
$ws.Range("D16").Value = "'0.001305"
$ws.Range("E16").Value = "'1.69%"
$ws.Range("D17").Value = "'0.006165"
$ws.Range("E17").Value = "'-1.07%"
$ws.Range("D19").Value = "'3.446"
$ws.Range("E19").Value = "'2.78%"
$ws.Range("E20").Value = "'-1.11%"
$ws.Range("D21").Value = "'8.317"
$ws.Range("E21").Value = "'-4.59%"
$ws.Range("D22").Value = "'0.1381"
$ws.Range("E22").Value = "'0.67%"
$ws.Range("D23").Value = "'0.2410"
$ws.Range("E23").Value = "'-1.51%"
$ws.Range("D24").Value = "'0.04416"
$ws.Range("E24").Value = "'-0.14%"
$ws.Range("D25").Value = "'0.001257"
$ws.Range("E25").Value = "'-0.40%"
$ws.Range("D26").Value = "'0.004321"
$ws.Range("E26").Value = "'-3.56%"
$ws.Range("D27").Value = "'0.0001181"
$ws.Range("E27").Value = "'2.52%"
$ws.Range("D39").Value = "'0.02774"
$ws.Range("E39").Value = "'10.84%"
$ws.Range("D40").Value = "'0.05401"
$ws.Range("E40").Value = "'2.68%"
$ws.Range("D41").Value = "'0.007660"
$ws.Range("E41").Value = "'2.30%"
$ws.Range("D42").Value = "'0.1418"
$ws.Range("E42").Value = "'1.01%"
$ws.Range("D43").Value = "'0.008960"
$ws.Range("E43").Value = "'-6.52%"
$ws.Range("D44").Value = "'0.002141"
$ws.Range("E44").Value = "'1.04%"
$ws.Range("D45").Value = "'0.01168"
$ws.Range("E45").Value = "'5.12%"
$ws.Range("D46").Value = "'0.00006559"
$ws.Range("E46").Value = "'-3.23%"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'-0.07%"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("E50").Value = "'-0.07%"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("E51").Value = "'-0.07%"
